$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 8.936014356752429
$ws.Range("C2").Value = 5.822997384038965
$ws.Range("D2").Value = 4.69325500866014
$ws.Range("E2").Value = 16.47667668692905
$ws.Range("F2").Value = 23.15625613623854
$ws.Range("I2").Value = 19.6611083466157
$ws.Range("K2").Value = 8.689240716366754
$ws.Range("O2").Value = 20.72629569191747
$ws.Range("B3").Value = 8.563451797084182
$ws.Range("C3").Value = 5.567647760220374
$ws.Range("D3").Value = 4.638264397795142
$ws.Range("E3").Value = 15.5442172489596
$ws.Range("F3").Value = 23.18289448938683
$ws.Range("I3").Value = 19.76091451388504
$ws.Range("K3").Value = 8.418152913604596
$ws.Range("O3").Value = 20.80098525690389
$ws.Range("B4").Value = 8.327146427361857
$ws.Range("C4").Value = 5.403620833193558
$ws.Range("D4").Value = 4.603726679817479
$ws.Range("E4").Value = 14.94693422056101
$ws.Range("F4").Value = 23.2065615617079
$ws.Range("I4").Value = 19.82664799666068
$ws.Range("K4").Value = 8.248216453697474
$ws.Range("O4").Value = 20.85225484366275
$ws.Range("B5").Value = 8.229090050531916
$ws.Range("C5").Value = 5.335010840441037
$ws.Range("D5").Value = 4.589465623672367
$ws.Range("E5").Value = 14.69758333982938
$ws.Range("F5").Value = 23.21803850255862
$ws.Range("I5").Value = 19.85455235888399
$ws.Range("K5").Value = 8.178193129835048
$ws.Range("O5").Value = 20.87450258663551
$ws.Range("B6").Value = 8.212705956464319
$ws.Range("C6").Value = 5.323513157659664
$ws.Range("D6").Value = 4.587086543558179
$ws.Range("E6").Value = 14.65582819829867
$ws.Range("F6").Value = 23.22005471693897
$ws.Range("I6").Value = 19.85925328745864
$ws.Range("K6").Value = 8.166522331598378
$ws.Range("O6").Value = 20.87827848372863
$ws.Range("B7").Value = 8.32583093508627
$ws.Range("C7").Value = 5.402702615761186
$ws.Range("D7").Value = 4.603535095592432
$ws.Range("E7").Value = 14.9435950957804
$ws.Range("F7").Value = 23.20670893313946
$ws.Range("I7").Value = 19.82701980295097
$ws.Range("K7").Value = 8.247275080429471
$ws.Range("O7").Value = 20.85254940557091
$ws.Range("B8").Value = 8.809199737762482
$ws.Range("C8").Value = 5.736485147693871
$ws.Range("D8").Value = 4.674459349130019
$ws.Range("E8").Value = 16.16044781658177
$ws.Range("F8").Value = 23.1639203540082
$ws.Range("I8").Value = 19.69459622013099
$ws.Range("K8").Value = 8.596553438774071
$ws.Range("O8").Value = 20.75092263419609
$ws.Range("B9").Value = 9.691932423644447
$ws.Range("C9").Value = 6.331634914449287
$ws.Range("D9").Value = 4.807074921354007
$ws.Range("E9").Value = 18.41992924820928
$ws.Range("F9").Value = 23.13824139689659
$ws.Range("I9").Value = 19.47033349738275
$ws.Range("K9").Value = 9.307098448307569
$ws.Range("O9").Value = 20.59479282944623
$ws.Range("B10").Value = 10.29481393372993
$ws.Range("C10").Value = 6.730563384020356
$ws.Range("D10").Value = 4.90013827490332
$ws.Range("E10").Value = 20.05566072055848
$ws.Range("F10").Value = 23.15510123277652
$ws.Range("I10").Value = 19.32729116686409
$ws.Range("K10").Value = 9.874003388772906
$ws.Range("O10").Value = 20.50670758719899
$ws.Range("B11").Value = 10.55818291148874
$ws.Range("C11").Value = 6.903392463923437
$ws.Range("D11").Value = 4.941440531158955
$ws.Range("E11").Value = 20.75765700656296
$ws.Range("F11").Value = 23.17055239400427
$ws.Range("I11").Value = 19.26696396349294
$ws.Range("K11").Value = 10.11979037569396
$ws.Range("O11").Value = 20.47247939673902
$ws.Range("B12").Value = 10.65628235904914
$ws.Range("C12").Value = 6.967573599696521
$ws.Range("D12").Value = 4.956925200635891
$ws.Range("E12").Value = 21.01746595111134
$ws.Range("F12").Value = 23.17752211719501
$ws.Range("I12").Value = 19.24480437313425
$ws.Range("K12").Value = 10.21109032953144
$ws.Range("O12").Value = 20.46036303645507
$ws.Range("B13").Value = 10.63522846765241
$ws.Range("C13").Value = 6.953807638903739
$ws.Range("D13").Value = 4.953597330320384
$ws.Range("E13").Value = 20.96177874503628
$ws.Range("F13").Value = 23.17597132618061
$ws.Range("I13").Value = 19.24954632428742
$ws.Range("K13").Value = 10.19150664301203
$ws.Range("O13").Value = 20.46293484758019
$ws.Range("B14").Value = 10.56628668345242
$ws.Range("C14").Value = 6.908698160177014
$ws.Range("D14").Value = 4.94271762232085
$ws.Range("E14").Value = 20.77915214270361
$ws.Range("F14").Value = 23.17110338304263
$ws.Range("I14").Value = 19.26512713510418
$ws.Range("K14").Value = 10.12733741125121
$ws.Range("O14").Value = 20.47146561132632
$ws.Range("B15").Value = 10.52384340062432
$ws.Range("C15").Value = 6.880901905912036
$ws.Range("D15").Value = 4.936033017639508
$ws.Range("E15").Value = 20.66650483313115
$ws.Range("F15").Value = 23.16826728518259
$ws.Range("I15").Value = 19.2747601278487
$ws.Range("K15").Value = 10.08779993137688
$ws.Range("O15").Value = 20.47680115365669
$ws.Range("B16").Value = 10.27737711967746
$ws.Range("C16").Value = 6.719092630259466
$ws.Range("D16").Value = 4.897417681597545
$ws.Range("E16").Value = 20.00893915501733
$ws.Range("F16").Value = 23.15424805405308
$ws.Range("I16").Value = 19.33132938920165
$ws.Range("K16").Value = 9.857694075128046
$ws.Range("O16").Value = 20.50906247384043
$ws.Range("B17").Value = 10.1233388617642
$ws.Range("C17").Value = 6.617596396422412
$ws.Range("D17").Value = 4.873458784065114
$ws.Range("E17").Value = 19.59478306624201
$ws.Range("F17").Value = 23.14764096899648
$ws.Range("I17").Value = 19.36724976252163
$ws.Range("K17").Value = 9.713406042121683
$ws.Range("O17").Value = 20.53035371951623
$ws.Range("B18").Value = 10.03371892567429
$ws.Range("C18").Value = 6.558405743543094
$ws.Range("D18").Value = 4.859581461735556
$ws.Range("E18").Value = 19.35260624838613
$ws.Range("F18").Value = 23.14457328220677
$ws.Range("I18").Value = 19.38835651086751
$ws.Range("K18").Value = 9.629278020698436
$ws.Range("O18").Value = 20.5431494936495
$ws.Range("B19").Value = 10.00320198771021
$ws.Range("C19").Value = 6.538225969186485
$ws.Range("D19").Value = 4.854866434248542
$ws.Range("E19").Value = 19.26992730444493
$ws.Range("F19").Value = 23.14366041652257
$ws.Range("I19").Value = 19.39557945816314
$ws.Range("K19").Value = 9.60059951377278
$ws.Range("O19").Value = 20.54757616313726
$ws.Range("B20").Value = 10.1398427091839
$ws.Range("C20").Value = 6.628485129812963
$ws.Range("D20").Value = 4.876019320084171
$ws.Range("E20").Value = 19.639280831394
$ws.Range("F20").Value = 23.14826848775436
$ws.Range("I20").Value = 19.3633797662319
$ws.Range("K20").Value = 9.728883728215871
$ws.Range("O20").Value = 20.52803031014985
$ws.Range("B21").Value = 10.58658135401009
$ws.Range("C21").Value = 6.921982410149992
$ws.Range("D21").Value = 4.945917534360897
$ws.Range("E21").Value = 20.83295717116981
$ws.Range("F21").Value = 23.17250286328947
$ws.Range("I21").Value = 19.26053205874475
$ws.Range("K21").Value = 10.14623386586777
$ws.Range("O21").Value = 20.46893694436724
$ws.Range("B22").Value = 10.869010066726
$ws.Range("C22").Value = 7.106416010773398
$ws.Range("D22").Value = 4.990690015025927
$ws.Range("E22").Value = 21.57800907555514
$ws.Range("F22").Value = 23.19486103952381
$ws.Range("I22").Value = 19.19730930848574
$ws.Range("K22").Value = 10.40864276766983
$ws.Range("O22").Value = 20.43524413119098
$ws.Range("B23").Value = 10.71916524691594
$ws.Range("C23").Value = 7.008662297335971
$ws.Range("D23").Value = 4.966879640748997
$ws.Range("E23").Value = 21.18356025048125
$ws.Range("F23").Value = 23.18233194334234
$ws.Range("I23").Value = 19.23068603544636
$ws.Range("K23").Value = 10.26954735900873
$ws.Range("O23").Value = 20.45277410366901
$ws.Range("B24").Value = 10.13238462021952
$ws.Range("C24").Value = 6.6235649454429
$ws.Range("D24").Value = 4.874862022114713
$ws.Range("E24").Value = 19.61917608125279
$ws.Range("F24").Value = 23.14798251000806
$ws.Range("I24").Value = 19.36512797238957
$ws.Range("K24").Value = 9.721889923330957
$ws.Range("O24").Value = 20.52907899455676
$ws.Range("B25").Value = 9.460773792683252
$ws.Range("C25").Value = 6.177232378861898
$ws.Range("D25").Value = 4.771936155038662
$ws.Range("E25").Value = 17.77999807161037
$ws.Range("F25").Value = 23.13892608037334
$ws.Range("I25").Value = 19.5271973436045
$ws.Range("K25").Value = 9.08786288563768
$ws.Range("O25").Value = 20.63237658057407
